$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.714.73"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.023.22"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'226.65"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'59.45"
$ws.Range("E7").Value = "  +3.83%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.326.66"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'14.50"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'20.88"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "'0.749"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "'5.20"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "2.023.54"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "37.747.20"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'6.01"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("D20").Value = "'69.44"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "'224.71"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'2.39"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").Value = "'2.20"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("D26").Value = "'164.60"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "'9.14"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("D29").Value = "'18.82"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "'1.27"
$ws.Range("E30").Value = "  -6.85%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "'4.41"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "'0.0599"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("D35").Value = "'4.46"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").Value = "'6.31"
$ws.Range("E36").Value = "  +5.80%  "
$ws.Range("D37").Value = "'2.23"
$ws.Range("E37").Value = "  -5.29%  "
$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "1.536.03"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("D41").Value = "'0.0215"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "'96.10"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'2.82"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'16.49"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "'0.0916"
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "'3.92"
$ws.Range("E47").Value = "  -5.10%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").Value = "'7.06"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "2.217.56"
$ws.Range("E51").Value = "  -1.41%  "
